# Error Calculations and Plots
# Two rows of source data were dropped (RM 232 and SC 92) and a number of
# individual A/B/C/D/E/F values were filled in / cleared elsewhere in the
# table. Delete from the bottom up first so the remaining row numbers don't
# need to be recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "SC 92" row (originally row 28) ...
$ws.Rows(28).Delete()
# ... and the "RM 232" row (originally row 26, still row 26 after the above).
$ws.Rows(26).Delete()

# Remaining per-cell corrections, addressed by the final (post-delete) row numbers.
$ws.Range("E2").Value = -7.2
$ws.Range("F3").ClearContents()
$ws.Range("F4").Value = 17.97
$ws.Range("F5").ClearContents()
$ws.Range("E6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("E12").Value = -5.3
$ws.Range("E14").ClearContents()
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("F19").ClearContents()
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F22").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").Value = 16.48
$ws.Range("E24").ClearContents()
$ws.Range("F25").Value = 16.6
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("E31").Value = -8.1
$ws.Range("B32").ClearContents()
$ws.Range("E33").Value = -10.7
